# Adds a "time_taken" column (F) to the panelapp worksheet, recording the
# per-row processing timestamp used when the panel data was generated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy the formatting from the existing header cell (E1,
# bold + border style) onto F1, then set its text.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps (F2:F36), one per data row, in file order.
$timestamps = @(
    "2021-10-05 13:39:14.532634",
    "2021-10-05 13:39:14.532647",
    "2021-10-05 13:39:14.532651",
    "2021-10-05 13:39:14.532655",
    "2021-10-05 13:39:14.532658",
    "2021-10-05 13:39:14.532662",
    "2021-10-05 13:39:14.532665",
    "2021-10-05 13:39:14.532668",
    "2021-10-05 13:39:14.532671",
    "2021-10-05 13:39:14.532674",
    "2021-10-05 13:39:14.532678",
    "2021-10-05 13:39:14.532681",
    "2021-10-05 13:39:14.532684",
    "2021-10-05 13:39:14.532687",
    "2021-10-05 13:39:14.532690",
    "2021-10-05 13:39:14.532693",
    "2021-10-05 13:39:14.532697",
    "2021-10-05 13:39:14.532700",
    "2021-10-05 13:39:14.532703",
    "2021-10-05 13:39:14.532706",
    "2021-10-05 13:39:14.532710",
    "2021-10-05 13:39:14.532713",
    "2021-10-05 13:39:14.532716",
    "2021-10-05 13:39:14.532719",
    "2021-10-05 13:39:14.532722",
    "2021-10-05 13:39:14.532726",
    "2021-10-05 13:39:14.532729",
    "2021-10-05 13:39:14.532732",
    "2021-10-05 13:39:14.532735",
    "2021-10-05 13:39:14.532738",
    "2021-10-05 13:39:14.532741",
    "2021-10-05 13:39:14.532744",
    "2021-10-05 13:39:14.532748",
    "2021-10-05 13:39:14.532751",
    "2021-10-05 13:39:14.532754"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
